# acpt: add scenarios for _Row.height, .height_rule
#
# Set an explicit row height / height-rule on the first (and only, for the
# single-cell tables) row of four of the tables in the document. Word's
# Row.Height is expressed in points, so twips-in-XML / 20 gives the value
# to assign (1 pt = 20 dxa/twips).
#   wdRowHeightAuto     = 0
#   wdRowHeightAtLeast  = 1
#   wdRowHeightExactly  = 2

$d = $word.ActiveDocument

# Table 2 (single 1440-dxa column): auto height, val=0
$row1 = $d.Tables.Item(2).Rows.Item(1)
$row1.HeightRule = 0
$row1.Height = 0

# Table 3 (single 2880-dxa column): at-least 2880 dxa (144 pt)
$row2 = $d.Tables.Item(3).Rows.Item(1)
$row2.HeightRule = 1
$row2.Height = 144

# Table 4 (3x1440-dxa columns, 2 rows): exactly 4320 dxa (216 pt) on row 1
$row3 = $d.Tables.Item(4).Rows.Item(1)
$row3.HeightRule = 2
$row3.Height = 216

# Table 5 (3x1440-dxa columns, 2 rows): exactly 5760 dxa (288 pt) on row 1
$row4 = $d.Tables.Item(5).Rows.Item(1)
$row4.HeightRule = 2
$row4.Height = 288
